# Updates the RTS.xlsx "Sheet1" address/ID listing:
#  - replaces the address (col A) and id (col N) text for existing rows 1-8
#  - appends 6 new rows (9-14) with the same A/N pattern (B:M stay blank)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new text content for column A (address) and column N (id), rows 1-8 ---
$addresses = @(
    "Rua Abílio Soares, 731 - Paraíso São Paulo/SP CEP:04005003",
    "Rua Desembargador do Vale, 836, ANEXO 830 - Perdizes São Paulo/SP CEP:05010040",
    "Rua Alves Guimarães, 1091, COZINHA 4 - Pinheiros São Paulo/SP CEP:05410-002",
    "Rua Doutor Augusto de Miranda, 549 - Vila Pompéia São Paulo/SP CEP:05026000",
    "Rua Guaipá, 1017,  - Vila Leopoldina São Paulo/SP CEP:05089-001",
    "Rua Coriolano, 301 - Vila Romana São Paulo/SP CEP:05047001",
    "Rua Catão, 479, NBURGER - Vila Romana São Paulo/SP CEP:05049000",
    "RUA CAMILO, 763, sem complemento - VILA ROMANA São Paulo/SP CEP:05045020",
    "Rua José Mariano Filho, 200,  - Jardim Oriental São Paulo/SP CEP:04347-180",
    "Rua Doutor Tomás Carvalhal, 626 - Paraíso São Paulo/SP CEP:04006001",
    "Rua Caraíbas, 964, IMOBILIARIA ESTEVAM - Perdizes São Paulo/SP CEP:05020000",
    "Rua Vergueiro, 4289,  - Vila Mariana São Paulo/SP CEP:04101-901",
    "Rua Borges Lagoa, 1050,  - Vila Clementino São Paulo/SP CEP:04038-002",
    "Rua Cardoso de Almeida, 587,  - Perdizes São Paulo/SP CEP:05013-000"
)

$ids = @(
    "72.791",
    "72.838",
    "72.859",
    "72.891",
    "72.903",
    "72.935",
    "72.937",
    "72.941",
    "72.956",
    "72.967",
    "72.970",
    "72.988",
    "72.990",
    "73.008"
)

# Force the id column to stay plain text (values look numeric, e.g.
# "72.791") instead of being parsed into a float.
$ws.Range("N1:N14").NumberFormat = "@"

for ($i = 0; $i -lt 14; $i++) {
    $row = $i + 1

    $ws.Cells.Item($row, 1).Value = $addresses[$i]
    $ws.Cells.Item($row, 14).Value = $ids[$i]
}
